$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 925
$ws.Range("I6").Value = 1000
$ws.Range("K6").Value = 3000
$ws.Range("M6").Value = -2888

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 199.125
$ws.Range("I8").Value = 199.125
$ws.Range("K8").Value = 597.375
$ws.Range("M8").Value = -458.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1624.75
$ws.Range("I31").Value = 1750
$ws.Range("J31").Value = 1499.5
$ws.Range("K31").Value = 5250
$ws.Range("L31").Value = 4498.5
$ws.Range("M31").Value = -5020
$ws.Range("N31").Value = -4958.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1613592
$ws.Range("J38").Value = 2066.6667
$ws.Range("L38").Value = 6200.000100000001
$ws.Range("N38").Value = -6944.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 992298.5
$ws.Range("I39").Value = 1133976.2
$ws.Range("J39").Value = 554
$ws.Range("K39").Value = 3401928.6
$ws.Range("L39").Value = 1662
$ws.Range("M39").Value = -3401632.6
$ws.Range("N39").Value = -2254

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 2395432.8
$ws.Range("I61").Value = 5714405.5
$ws.Range("J61").Value = 24737.857
$ws.Range("K61").Value = 17143216.5
$ws.Range("L61").Value = 74213.571
$ws.Range("M61").Value = -17143044.5
$ws.Range("N61").Value = -74557.571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1806.909
$ws.Range("I70").Value = 1867.6
$ws.Range("J70").Value = 1200
$ws.Range("K70").Value = 5602.799999999999
$ws.Range("L70").Value = 3600
$ws.Range("M70").Value = -5332.799999999999
$ws.Range("N70").Value = -4140

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1806.909
$ws.Range("I73").Value = 1867.6
$ws.Range("J73").Value = 1200
$ws.Range("K73").Value = 5602.799999999999
$ws.Range("L73").Value = 3600
$ws.Range("M73").Value = -4666.799999999999
$ws.Range("N73").Value = -5472

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1645644.6
$ws.Range("J112").Value = 1737036
$ws.Range("L112").Value = 5211108
$ws.Range("N112").Value = -5213324

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1949.2059
$ws.Range("I137").Value = 1283.8966
$ws.Range("J137").Value = 5808
$ws.Range("K137").Value = 3851.6898
$ws.Range("L137").Value = 17424
$ws.Range("M137").Value = -1301.6898
$ws.Range("N137").Value = -22524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24365.393
$ws.Range("I32").Value = 4075.194
$ws.Range("J32").Value = 218571.58
$ws.Range("K32").Value = 4075.194
$ws.Range("L32").Value = 218571.58
$ws.Range("M32").Value = -3788.194
$ws.Range("N32").Value = -219145.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 36940.57
$ws.Range("I45").Value = 50942.9
$ws.Range("J45").Value = 1934.75
$ws.Range("K45").Value = 50942.9
$ws.Range("L45").Value = 1934.75
$ws.Range("M45").Value = -50565.9
$ws.Range("N45").Value = -2688.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2635.3928
$ws.Range("I61").Value = 2040.1
$ws.Range("K61").Value = 2040.1
$ws.Range("M61").Value = -1828.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2828
$ws.Range("I74").Value = 875
$ws.Range("J74").Value = 4130
$ws.Range("K74").Value = 875
$ws.Range("L74").Value = 4130
$ws.Range("M74").Value = -1
$ws.Range("N74").Value = -5878

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2828
$ws.Range("I77").Value = 875
$ws.Range("J77").Value = 4130
$ws.Range("K77").Value = 4375
$ws.Range("L77").Value = 20650
$ws.Range("M77").Value = -7
$ws.Range("N77").Value = -29386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2635.3928
$ws.Range("I136").Value = 2040.1
$ws.Range("K136").Value = 6120.299999999999
$ws.Range("M136").Value = -3570.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2814
$ws.Range("I134").Value = 2964.9375
$ws.Range("J134").Value = 2442.4614
$ws.Range("K134").Value = 8894.8125
$ws.Range("L134").Value = 7327.3842
$ws.Range("M134").Value = -6359.8125
$ws.Range("N134").Value = -12397.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 1006
$ws.Range("J11").Value = 1006
$ws.Range("L11").Value = 1006
$ws.Range("N11").Value = -1286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29462.508
$ws.Range("I31").Value = 1069.6
$ws.Range("J31").Value = 44810.027
$ws.Range("K31").Value = 1069.6
$ws.Range("L31").Value = 44810.027
$ws.Range("M31").Value = -774.5999999999999
$ws.Range("N31").Value = -45400.027

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 29462.508
$ws.Range("I34").Value = 1069.6
$ws.Range("J34").Value = 44810.027
$ws.Range("K34").Value = 1069.6
$ws.Range("L34").Value = 44810.027
$ws.Range("M34").Value = -867.5999999999999
$ws.Range("N34").Value = -45214.027

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4963.7
$ws.Range("I58").Value = 1182.7812
$ws.Range("J58").Value = 20087.375
$ws.Range("K58").Value = 1182.7812
$ws.Range("L58").Value = 20087.375
$ws.Range("M58").Value = -979.7811999999999
$ws.Range("N58").Value = -20493.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 18375
$ws.Range("J95").Value = 18375
$ws.Range("L95").Value = 18375
$ws.Range("N95").Value = -23867

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 39003
$ws.Range("I99").Value = 4980
$ws.Range("K99").Value = 4980
$ws.Range("M99").Value = -3482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 39003
$ws.Range("I126").Value = 4980
$ws.Range("K126").Value = 14940
$ws.Range("M126").Value = -12470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H130").Value = 48240
$ws.Range("J130").Value = 48240
$ws.Range("L130").Value = 48240
$ws.Range("N130").Value = -58280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 23440038
$ws.Range("I132").Value = 22224656
$ws.Range("J132").Value = 26318574
$ws.Range("K132").Value = 66673968
$ws.Range("L132").Value = 78955722
$ws.Range("M132").Value = -66671438
$ws.Range("N132").Value = -78960782

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1282.8235
$ws.Range("I134").Value = 1224.8889
$ws.Range("J134").Value = 1348
$ws.Range("K134").Value = 3674.6667
$ws.Range("L134").Value = 4044
$ws.Range("M134").Value = -1139.6667
$ws.Range("N134").Value = -9114

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4963.7
$ws.Range("I136").Value = 1182.7812
$ws.Range("J136").Value = 20087.375
$ws.Range("K136").Value = 3548.3436
$ws.Range("L136").Value = 60262.125
$ws.Range("M136").Value = -998.3435999999997
$ws.Range("N136").Value = -65362.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 678.8889
$ws.Range("J34").Value = 958.3333
$ws.Range("L34").Value = 2874.9999
$ws.Range("N34").Value = -3042.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10597.5
$ws.Range("J122").Value = 11732.777
$ws.Range("L122").Value = 105594.993
$ws.Range("N122").Value = -110494.993

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6948.828
$ws.Range("I131").Value = 2126.6667
$ws.Range("J131").Value = 7099.521
$ws.Range("K131").Value = 6380.000100000001
$ws.Range("L131").Value = 21298.563
$ws.Range("M131").Value = -1340.000100000001
$ws.Range("N131").Value = -31378.563

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2154.318
$ws.Range("J132").Value = 2259.75
$ws.Range("L132").Value = 20337.75
$ws.Range("N132").Value = -25397.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 166840030
$ws.Range("I80").Value = 200207600
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 200207600
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -200206602
$ws.Range("N80").Value = -4196

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 166840030
$ws.Range("I83").Value = 200207600
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 1001038000
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -1001033008
$ws.Range("N83").Value = -20984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 33000
$ws.Range("J93").Value = 33000
$ws.Range("L93").Value = 33000
$ws.Range("N93").Value = -36744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6980
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 6980
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 6980
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -7384

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 6000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 6000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -7082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6980
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 6980
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6980
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -11320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2667.6445
$ws.Range("I132").Value = 2706.342
$ws.Range("J132").Value = 2457.5715
$ws.Range("K132").Value = 8119.026
$ws.Range("L132").Value = 7372.7145
$ws.Range("M132").Value = -5589.026
$ws.Range("N132").Value = -12432.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1951.5
$ws.Range("I136").Value = 1727.4546
$ws.Range("J136").Value = 2303.5715
$ws.Range("K136").Value = 5182.3638
$ws.Range("L136").Value = 6910.7145
$ws.Range("M136").Value = -2632.3638
$ws.Range("N136").Value = -12010.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2017.4445
$ws.Range("I122").Value = 1850
$ws.Range("K122").Value = 5550
$ws.Range("M122").Value = -3100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1888.1613
$ws.Range("I132").Value = 1838.9791
$ws.Range("J132").Value = 2056.7856
$ws.Range("K132").Value = 5516.9373
$ws.Range("L132").Value = 6170.3568
$ws.Range("M132").Value = -2986.9373
$ws.Range("N132").Value = -11230.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 35740.875
$ws.Range("J135").Value = 35740.875
$ws.Range("L135").Value = 35740.875
$ws.Range("N135").Value = -45880.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1200.6522
$ws.Range("I136").Value = 807.8570999999999
$ws.Range("J136").Value = 1811.6666
$ws.Range("K136").Value = 2423.5713
$ws.Range("L136").Value = 5434.9998
$ws.Range("M136").Value = 126.4287000000004
$ws.Range("N136").Value = -10534.9998
